# "break out stock.yaml completed"
#
# 1) Sheet "10per change": columns D2:D10 were stored as text (inlineStr);
#    convert them to plain numeric values (same digits), then append 11 new
#    data rows (11-21) pulled from a later screener run.
# 2) Sheet "DND 3 V 0.3": column D2 was stored as text; convert it to a
#    plain numeric value the same way.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as genuine text (t="s"/inlineStr),
# even when it looks like a number (e.g. a BSE code). Excel normally
# auto-detects numeric-looking strings and stores them as numbers, so we
# force text via a leading apostrophe and then strip the quote-prefix
# style Excel applies for that, leaving a plain text cell behind.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.ClearFormats()
}

# -----------------------------------------------------------------
# 1) "10per change" sheet
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("10per change")

# D2:D10 switch from text to numeric (values unchanged).
$ws1.Range("D2").Value = 500510
$ws1.Range("D3").Value = 543287
$ws1.Range("D4").Value = 542066
$ws1.Range("D5").Value = 500112
$ws1.Range("D6").Value = 533096
$ws1.Range("D7").Value = 500093
$ws1.Range("D8").Value = 543396
$ws1.Range("D9").Value = 532898
$ws1.Range("D10").Value = 532155

# New rows 11-21 (bsecode column D stays text, like the original rows did
# before this commit).
$newRows = @(
    @(11, 1,  "HDFCAMC",    "HDFC Asset Management Company Ltd",             "541729", 1.29, 3753.35, 373559,   "06/06/2024 09:26:40"),
    @(12, 2,  "LT",         "Larsen & Toubro Limited",                       "500510", 2.07, 3479.5,  7241361,  "06/06/2024 09:26:40"),
    @(13, 3,  "ADANIENT",   "Adani Enterprises Limited",                     "512599", 2.22, 3184.4,  5307540,  "06/06/2024 09:26:40"),
    @(14, 4,  "ADANIGREEN", "Adani Green Energy Ltd",                        "541450", 1.98, 1865,    1722796,  "06/06/2024 09:26:40"),
    @(15, 5,  "ADANIPORTS", "Adani Ports And Special Economic Zone Limited", "532921", 0.07000000000000001, 1355.6, 10277842, "06/06/2024 09:26:40"),
    @(16, 6,  "SBIN",       "State Bank Of India",                           "500112", 3.38, 816.45,  36455863, "06/06/2024 09:26:40"),
    @(17, 7,  "ADANIPOWER", "Adani Power Limited",                           "533096", 4.04, 756,     15321524, "06/06/2024 09:26:40"),
    @(18, 8,  "COALINDIA",  "Coal India Limited",                            "533278", 2.43, 471.75,  18526644, "06/06/2024 09:26:40"),
    @(19, 9,  "POWERGRID",  "Power Grid Corporation Of India Limited",       "532898", 0.27, 299.6,   30068761, "06/06/2024 09:26:40"),
    @(20, 10, "BANKBARODA", "Bank Of Baroda",                                "532134", 3.37, 268.75,  29153860, "06/06/2024 09:26:40"),
    @(21, 11, "GAIL",       "Gail (india) Limited",                          "532155", 6,    206.85,  37330496, "06/06/2024 09:26:40")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    Set-TextValue $ws1.Cells.Item($r, 4) $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
    $ws1.Cells.Item($r, 6).Value = $row[6]
    $ws1.Cells.Item($r, 7).Value = $row[7]
    $ws1.Cells.Item($r, 8).Value = $row[8]
}

# -----------------------------------------------------------------
# 2) "DND 3 V 0.3" sheet
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("DND 3 V 0.3")
$ws3.Range("D2").Value = 500331
